$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "dimension 9 x 9" -> "size 9 x 9"  (Representation section)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("dimension 9 x 9", $true, $false, $false, $false, $false, $true, 1, $false, "size 9 x 9", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: "average deviation" -> "standard deviation" (Evaluation section, 2 occurrences)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("average deviation", $true, $false, $false, $false, $false, $true, 1, $false, "standard deviation", 2) | Out-Null
$d.Content.Find.Execute("average deviation", $true, $false, $false, $false, $false, $true, 1, $false, "standard deviation", 2) | Out-Null

# Change 2b: "the estimated value of 45" -> "the value of 45"
$d.Content.Find.Execute("the estimated value of 45", $true, $false, $false, $false, $false, $true, 1, $false, "the value of 45", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3: recombination text tweaks
# ---------------------------------------------------------------------
$d.Content.Find.Execute("a single slice (Slice after", $true, $false, $false, $false, $false, $true, 1, $false, "a single slice (only one slice after", 2) | Out-Null
$d.Content.Find.Execute("chance of a slice after every 3x3 subgrid)", $true, $false, $false, $false, $false, $true, 1, $false, "chance of a slice between every 3x3 subgrid)", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 4: mutation text tweak
# ---------------------------------------------------------------------
$d.Content.Find.Execute("happens with a probability or 1/9, since we have 9 3x3 subgrids in a Sudoku", $true, $false, $false, $false, $false, $true, 1, $false, "is performed on one random 3x3 subgrid of a Sudoku", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 5: append formula explanation to the Selection paragraph
# ---------------------------------------------------------------------
$d.Content.Find.Execute("get selected).", $true, $false, $false, $false, $false, $true, 1, $false, "get selected). The formula used for this calculation is: (size- ranking)/size. With size being the number of parents and children combined, and ranking being the position of the Sudoku in the list, sorted by fitness.", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 6: drop the stale lastRenderedPageBreak marker before "Additional Info:"
# A find/replace that round-trips the exact same text forces the run to be
# rewritten, which clears the stale cached page-break marker (just like a
# real edit in Word would invalidate it).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Additional Info:", $true, $false, $false, $false, $false, $true, 1, $false, "Additional Infox:", 2) | Out-Null
$d.Content.Find.Execute("Additional Infox:", $true, $false, $false, $false, $false, $true, 1, $false, "Additional Info:", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 7: "combined to a new initial population" -> "combined to form a new initial population"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("combined to a new initial population", $true, $false, $false, $false, $false, $true, 1, $false, "combined to form a new initial population", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 8: add two new paragraphs after the "solveMultiple" paragraph:
#   - an empty paragraph (Listenabsatz style, no numbering)
#   - a numbered paragraph (Listenabsatz, same list as the others) with new text
# ---------------------------------------------------------------------
$solveMultiplePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*with this population.*") {
        $solveMultiplePara = $p
        break
    }
}

# The paragraph right after solveMultiplePara is the existing empty
# "Listenabsatz" paragraph (no numbering) that precedes the final
# "To run the project..." item. Duplicate its (paragraph-mark) formatting
# to create a matching blank paragraph right after solveMultiplePara.
$emptyPara = $solveMultiplePara.Next()
$emptyCopy = $emptyPara.Range.FormattedText
$insertPoint1 = $d.Range($solveMultiplePara.Range.End, $solveMultiplePara.Range.End)
$insertPoint1.FormattedText = $emptyCopy

# Newly created blank paragraph now sits between solveMultiplePara and the
# original emptyPara.
$newEmptyPara = $solveMultiplePara.Next()

# Duplicate solveMultiplePara itself (keeps the Listenabsatz + numPr
# formatting, incl. numId) right after the new blank paragraph, then
# overwrite its text with the new sentence.
$numberedCopy = $solveMultiplePara.Range.FormattedText
$insertPoint2 = $d.Range($newEmptyPara.Range.End, $newEmptyPara.Range.End)
$insertPoint2.FormattedText = $numberedCopy

$newNumberedPara = $newEmptyPara.Next()
$newText = "The algorithm stops when either a correct solution for a Sudoku is found or after a certain number of iterations (in case no correct solution could be found). If the latter is the case the best Sudoku calculated so far is returned. "
$r = $newNumberedPara.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = $newText
